$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (and E) to Text format while writing so that numeric-looking
# strings (e.g. "545.08", "11.00", "0.0946") keep their exact text representation
# instead of being auto-converted to floating point numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '58.463.66'
$ws.Range('E2').Value = '  -2.79%  '
$ws.Range('D3').Value = '2.280.29'
$ws.Range('E3').Value = '  -5.65%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '545.08'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').Value = '130.73'
$ws.Range('E6').Value = '  -4.65%  '
$ws.Range('E8').Value = '  -3.26%  '
$ws.Range('E9').Value = '  -3.71%  '
$ws.Range('E10').Value = '  -2.54%  '
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E12').Value = '  -5.51%  '
$ws.Range('D13').Value = '23.60'
$ws.Range('E13').Value = '  -5.18%  '
$ws.Range('D14').Value = '2.685.18'
$ws.Range('E14').Value = '  -5.73%  '
$ws.Range('D15').Value = '58.422.78'
$ws.Range('E15').Value = '  -2.70%  '
$ws.Range('E16').Value = '  -3.41%  '
$ws.Range('D17').Value = '2.279.76'
$ws.Range('E17').Value = '  -6.59%  '
$ws.Range('E18').Value = '  -6.06%  '
$ws.Range('D19').Value = '4.30'
$ws.Range('E19').Value = '  -3.96%  '
$ws.Range('D20').Value = '314.28'
$ws.Range('E20').Value = '  -4.20%  '
$ws.Range('D21').Value = '6.44'
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '63.02'
$ws.Range('E23').Value = '  -3.74%  '
$ws.Range('E24').Value = '  -4.31%  '
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -1.91%  '
$ws.Range('E26').Value = '  -6.82%  '
$ws.Range('E27').Value = '  -5.93%  '
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('D29').Value = '170.52'
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').Value = '0.0₃0722'
$ws.Range('E30').Value = '  -6.60%  '
$ws.Range('D31').Value = '1.08'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').Value = '5.75'
$ws.Range('E32').Value = '  -5.70%  '
$ws.Range('D33').Value = '0.384'
$ws.Range('E33').Value = '  -5.04%  '
$ws.Range('D35').Value = '17.77'
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -5.00%  '
$ws.Range('E38').Value = '  -6.05%  '
$ws.Range('D39').Value = '37.94'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('E40').Value = '  -4.99%  '
$ws.Range('D41').Value = '300.32'
$ws.Range('E41').Value = '  -8.83%  '
$ws.Range('D42').Value = '140.82'
$ws.Range('E42').Value = '  -2.50%  '
$ws.Range('D43').Value = '3.45'
$ws.Range('E43').Value = '  -5.40%  '
$ws.Range('D44').Value = '0.0946'
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').Value = '0.0498'
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('D46').Value = '0.552'
$ws.Range('E46').Value = '  -4.30%  '
$ws.Range('D47').Value = '18.39'
$ws.Range('E47').Value = '  -8.32%  '
$ws.Range('E48').Value = '  -4.15%  '
$ws.Range('E49').Value = '  -5.72%  '
$ws.Range('D50').Value = '11.00'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '4.61'
$ws.Range('E51').Value = '  -0.88%  '

# Restore default (General) formatting/style so the cells match the original
# unstyled appearance (no explicit style index) instead of keeping a Text format.
$dataRange.ClearFormats()
